$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -65989.8900122636
$ws.Range("D4").Value = -113313.936254694
$ws.Range("E4").Value = 297905.677930427
$ws.Range("F4").Value = 346789.301654904
$ws.Range("C5").Value = -4582526.53912444
$ws.Range("D5").Value = -7148990.8544689
$ws.Range("E5").Value = 5113806.45582355
$ws.Range("F5").Value = 7680270.77116801
$ws.Range("C6").Value = -164694572.009674
$ws.Range("D6").Value = -252083589.40136
$ws.Range("E6").Value = 165468994.911911
$ws.Range("F6").Value = 252858012.303597
$ws.Range("C7").Value = -3135989273.58038
$ws.Range("D7").Value = -4796279150.61706
$ws.Range("E7").Value = 3136735197.84024
$ws.Range("F7").Value = 4797025074.87693
$ws.Range("C8").Value = -36931461439.3753
$ws.Range("D8").Value = -56481956677.0706
$ws.Range("E8").Value = 36932067816.5906
$ws.Range("F8").Value = 56482563054.2859
$ws.Range("C9").Value = -217695560157.655
$ws.Range("D9").Value = -332936716972.214
$ws.Range("E9").Value = 217695894348.896
$ws.Range("F9").Value = 332937051163.455
$ws.Range("C10").Value = -2018085365820.55
$ws.Range("D10").Value = -3086395317172.28
$ws.Range("E10").Value = 2018085707430.0
$ws.Range("F10").Value = 3086395658781.74
$ws.Range("C11").Value = -17082522575490.2
$ws.Range("D11").Value = -26125463876753.8
$ws.Range("E11").Value = 17082522925659.6
$ws.Range("F11").Value = 26125464226923.3
$ws.Range("C12").Value = -136405555098887.0
$ws.Range("D12").Value = -208614294244524.0
$ws.Range("E12").Value = 136405555390864.0
$ws.Range("F12").Value = 208614294536502.0
$ws.Range("C13").Value = -1129927187562564.0
$ws.Range("D13").Value = -1728074509358073.0
$ws.Range("E13").Value = 1129927187769806.0
$ws.Range("F13").Value = 1728074509565316.0
$ws.Range("C14").Value = -7618687840764581.0
$ws.Range("D14").Value = -11651777563116704.0
$ws.Range("E14").Value = 7618687840876191.0
$ws.Range("F14").Value = 11651777563228312.0
$ws.Range("C15").Value = -59577841479199784.0
$ws.Range("D15").Value = -91116445654895008.0
$ws.Range("E15").Value = 59577841479265032.0
$ws.Range("F15").Value = 91116445654960256.0
$ws.Range("C16").Value = -3510747451053496832.0
$ws.Range("D16").Value = -5369224889485145088.0
$ws.Range("E16").Value = 3510747451053715968.0
$ws.Range("F16").Value = 5369224889485364224.0
$ws.Range("C17").Value = -92303315148741033984.0
$ws.Range("D17").Value = -141165738632065089536.0
$ws.Range("E17").Value = 92303315148741558272.0
$ws.Range("F17").Value = 141165738632065613824.0
$ws.Range("C18").Value = -1424255128650531471360.0
$ws.Range("D18").Value = -2178210250763688476672.0
$ws.Range("E18").Value = 1424255128650531995648.0
$ws.Range("F18").Value = 2178210250763689000960.0
$ws.Range("C19").Value = -19882790989398091497472.0
$ws.Range("D19").Value = -30408104752926959861760.0
$ws.Range("E19").Value = 19882790989398091497472.0
$ws.Range("F19").Value = 30408104752926959861760.0
$ws.Range("C20").Value = -123027017129205574926336.0
$ws.Range("D20").Value = -188153585998052791943168.0
$ws.Range("E20").Value = 123027017129205574926336.0
$ws.Range("F20").Value = 188153585998052791943168.0
$ws.Range("C21").Value = -211900720638101584609280.0
$ws.Range("D21").Value = -324074186255838490918912.0
$ws.Range("E21").Value = 211900720638101584609280.0
$ws.Range("F21").Value = 324074186255838490918912.0
$ws.Range("C22").Value = -1955540409863661418446848.0
$ws.Range("D22").Value = -2990740971095229982572544.0
$ws.Range("E22").Value = 1955540409863661418446848.0
$ws.Range("F22").Value = 2990740971095229982572544.0
$ws.Range("C23").Value = -7492950052875102138138624.0
$ws.Range("D23").Value = -11459478210969879659413504.0
$ws.Range("E23").Value = 7492950052875102138138624.0
$ws.Range("F23").Value = 11459478210969879659413504.0
$ws.Range("C24").Value = -50900814014913095562952704.0
$ws.Range("D24").Value = -77846077313795995121745920.0
$ws.Range("E24").Value = 50900814014913095562952704.0
$ws.Range("F24").Value = 77846077313795995121745920.0
$ws.Range("C25").Value = -267688491383489764603723776.0
$ws.Range("D25").Value = -409394218924420743001276416.0
$ws.Range("E25").Value = 267688491383489764603723776.0
$ws.Range("F25").Value = 409394218924420743001276416.0
$ws.Range("C26").Value = -3798298332012455521895317504.0
$ws.Range("D26").Value = -5808996011892340809221013504.0
$ws.Range("E26").Value = 3798298332012455521895317504.0
$ws.Range("F26").Value = 5808996011892340809221013504.0
$ws.Range("C27").Value = -56496680761477166749301342208.0
$ws.Range("D27").Value = -86404216978578144373728346112.0
$ws.Range("E27").Value = 56496680761477166749301342208.0
$ws.Range("F27").Value = 86404216978578144373728346112.0
$ws.Range("C28").Value = -3510788078359082693435932540928.0
$ws.Range("D28").Value = -5369287023587121941899825381376.0
$ws.Range("E28").Value = 3510788078359082693435932540928.0
$ws.Range("F28").Value = 5369287023587121941899825381376.0
$ws.Range("C29").Value = -68646344016596019595657714073600.0
$ws.Range("D29").Value = -104985523454688013534209683488768.0
$ws.Range("E29").Value = 68646344016596019595657714073600.0
$ws.Range("F29").Value = 104985523454688013534209683488768.0
$ws.Range("C30").Value = -971516363818411705926155835015168.0
$ws.Range("D30").Value = -1485806060926021133579936396738560.0
$ws.Range("E30").Value = 971516363818411705926155835015168.0
$ws.Range("F30").Value = 1485806060926021133579936396738560.0
$ws.Range("C31").Value = -9652176975403372029600599424106496.0
$ws.Range("D31").Value = -14761730821309643375838662508937216.0
$ws.Range("E31").Value = 9652176975403372029600599424106496.0
$ws.Range("F31").Value = 14761730821309643375838662508937216.0
$ws.Range("C32").Value = -43826591180547271265208620516638720.0
$ws.Range("D32").Value = -67026987121295013006532546012381184.0
$ws.Range("E32").Value = 43826591180547271265208620516638720.0
$ws.Range("F32").Value = 67026987121295013006532546012381184.0
$ws.Range("C33").Value = -177132585059878873352991694246641664.0
$ws.Range("D33").Value = -270900911473123119977607167968542720.0
$ws.Range("E33").Value = 177132585059878873352991694246641664.0
$ws.Range("F33").Value = 270900911473123119977607167968542720.0
$ws.Range("C34").Value = -1148234671426637444276653644939526144.0
$ws.Range("D34").Value = -1756073389711816241272403176280031232.0
$ws.Range("E34").Value = 1148234671426637444276653644939526144.0
$ws.Range("F34").Value = 1756073389711816241272403176280031232.0
$ws.Range("C35").Value = -3395371083304086024239377598619058176.0
$ws.Range("D35").Value = -5192771962005889163314238944083705856.0
$ws.Range("E35").Value = 3395371083304086024239377598619058176.0
$ws.Range("F35").Value = 5192771962005889163314238944083705856.0
$ws.Range("C36").Value = -31738718042627063525742703736636047360.0
$ws.Range("D36").Value = -48540180474584037221489227309976125440.0
$ws.Range("E36").Value = 31738718042627063525742703736636047360.0
$ws.Range("F36").Value = 48540180474584037221489227309976125440.0
$ws.Range("C37").Value = -167783837345194140686634776984094244864.0
$ws.Range("D37").Value = -256603235660486824210286045644785188864.0
$ws.Range("E37").Value = 167783837345194140686634776984094244864.0
$ws.Range("F37").Value = 256603235660486824210286045644785188864.0
$ws.Range("C38").Value = -708994583697721355540217614264812175360.0
$ws.Range("D38").Value = -1084313644992493545144236543218624233472.0
$ws.Range("E38").Value = 708994583697721355540217614264812175360.0
$ws.Range("F38").Value = 1084313644992493545144236543218624233472.0
$ws.Range("C39").Value = -13489612197994482408427288678821391237120.0
$ws.Range("D39").Value = -20630581542183945020318019647869133258752.0
$ws.Range("E39").Value = 13489612197994482408427288678821391237120.0
$ws.Range("F39").Value = 20630581542183945020318019647869133258752.0
$ws.Range("C40").Value = -882558621779350216999103158358090379689984.0
$ws.Range("D40").Value = -1349756934827475986212792713356669238640640.0
$ws.Range("E40").Value = 882558621779350216999103158358090379689984.0
$ws.Range("F40").Value = 1349756934827475986212792713356669238640640.0
